# Add four new brand names to the restriction list, each inserted in its
# correct alphabetically-sorted position (matching how the list was
# originally built/sorted), pushing the rows below them down by one.
#
# Insertions (in top-to-bottom sheet order):
#   row 12  -> "bioderma"      (before "boss")
#   row 70  -> "juliette"      (after "juliette has a gun", before "kate")
#   row 98  -> "l'occittane"   (after "l'occitanne", before "loreal")
#   row 141 -> "thierry"       (before "thierry mugler")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-Brand($rowNum, $text) {
    $ws.Rows.Item($rowNum).Insert()
    $ws.Range("A" + $rowNum).Value = $text
    # Match formatting (style) of the row immediately below, which holds
    # the row that used to occupy this slot before the insert.
    $ws.Cells.Item($rowNum + 1, 1).Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)
}

Insert-Brand 12 "bioderma"
Insert-Brand 70 "juliette"
Insert-Brand 98 "l'occittane"
Insert-Brand 141 "thierry"

$excel.CutCopyMode = 0

$ws.Range("C216").Select()
